$d = $word.ActiveDocument

function Find-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. "Ported a WiFi Microcontroller BootRom ..." bullet -> new text
#    describing the encryption-API work.
# ---------------------------------------------------------------------
$p12 = Find-ParagraphByText $d "Ported a"
if ($p12 -ne $null) {
    $r12 = $p12.Range
    $r12.End = $r12.End - 1   # exclude the paragraph mark
    $r12.Text = "Adapted encryption API to run on dedicated cryptography hardware for new board release."
}

# ---------------------------------------------------------------------
# 2. "Composed memory map of code ..." bullet -> new text describing the
#    board startup code / linker script work. This paragraph is also
#    where the _GoBack bookmark now lives (right before "interrupts").
# ---------------------------------------------------------------------
$p13 = Find-ParagraphByText $d "Composed memory map"
if ($p13 -ne $null) {
    $r13 = $p13.Range
    $r13.End = $r13.End - 1   # exclude the paragraph mark
    $prefix13 = "Wrote board startup code to initialize cache, setup "
    $r13.Text = $prefix13 + "interrupts, and implement timers. Composed linker script to assemble code, data, stack, and heap for boot code."
}

# ---------------------------------------------------------------------
# 3. Move the _GoBack bookmark from its old spot (end of the "Currently
#    integrating bluetooth ... patch device settings." bullet) to the
#    new spot inside the rewritten paragraph, right before "interrupts".
# ---------------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
if ($oldBm -ne $null) {
    $oldBm.Delete()
}

$p13after = Find-ParagraphByText $d "Wrote board startup code"
if ($p13after -ne $null) {
    $bmPos = $p13after.Range.Start + $prefix13.Length
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
